# -------------------------------------------------------------------------
# Apply "Add files via upload" edit to FOCUSEDGROWTH_holdings.xlsx:
#   1) Bump the "as of" date in the confidentiality footnote from
#      2021-05-11 -> 2021-05-12
#   2) Refresh the Weight (col D) and Percent Change (col E) figures for
#      every holding row (2-56) with the latest model values
# -------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected; unprotect so the cells can be edited, then
# restore protection afterwards.
$ws.Unprotect()

# --- 1) Footnote date -------------------------------------------------
$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-12 for illustrative purposes only and are subject to change."

# --- 2) Weight / Percent Change refresh --------------------------------
$ws.Range("D2").Value = 0.01533966108397206
$ws.Range("E2").Value = -0.02102713711675008
$ws.Range("D3").Value = 0.05041058040242723
$ws.Range("E3").Value = -0.02232382417623313
$ws.Range("D4").Value = 0.01451064657929423
$ws.Range("E4").Value = -0.01987068965517258
$ws.Range("D5").Value = 0.009924781893113833
$ws.Range("E5").Value = -0.0181497353163601
$ws.Range("D6").Value = 0.01548667519942306
$ws.Range("E6").Value = -0.009749399244765034
$ws.Range("D7").Value = 0.02013517959486709
$ws.Range("E7").Value = -0.01058471421271645
$ws.Range("D8").Value = 0.00462842716461665
$ws.Range("E8").Value = 0.003263491462895729
$ws.Range("D9").Value = 0.006821379901857706
$ws.Range("E9").Value = -0.009077405512460834
$ws.Range("D10").Value = 0.0143613182658974
$ws.Range("E10").Value = -0.01936305732484078
$ws.Range("D11").Value = 0.008394283954340859
$ws.Range("E11").Value = -0.009164741822517031
$ws.Range("D12").Value = 0.01562271251110394
$ws.Range("E12").Value = -0.0165745856353593
$ws.Range("D13").Value = 0.002832703377483775
$ws.Range("E13").Value = -0.04250386398763517
$ws.Range("D14").Value = 0.005844600731043748
$ws.Range("E14").Value = -0.02825191288993512
$ws.Range("D15").Value = 0.01437823692925817
$ws.Range("E15").Value = -0.00687523653336708
$ws.Range("D16").Value = 0.01053632433677353
$ws.Range("E16").Value = -0.008548150126886478
$ws.Range("D17").Value = 0.02080526499196739
$ws.Range("E17").Value = -0.01001082251082241
$ws.Range("D18").Value = 0.00844031772969586
$ws.Range("E18").Value = 0.004201680672268893
$ws.Range("D19").Value = 0.01690039996076697
$ws.Range("E19").Value = -0.004026527711984862
$ws.Range("D20").Value = 0.01218081216085152
$ws.Range("E20").Value = -0.007958921694480092
$ws.Range("D21").Value = 0.007465602573266027
$ws.Range("E21").Value = -0.008989460632362056
$ws.Range("D22").Value = 0.01477559097114721
$ws.Range("E22").Value = -0.0117340286831813
$ws.Range("D23").Value = 0.01995836236262716
$ws.Range("E23").Value = -0.01937640140457753
$ws.Range("D24").Value = 0.01020551912229987
$ws.Range("E24").Value = -0.02194970230160831
$ws.Range("D25").Value = 0.02012414024520655
$ws.Range("E25").Value = -0.02715151515151515
$ws.Range("D26").Value = 0.01407573373020319
$ws.Range("E26").Value = -0.01693649979004352
$ws.Range("D27").Value = 0.02037973402644054
$ws.Range("E27").Value = -0.05867359856860721
$ws.Range("D28").Value = 0.05512607144651965
$ws.Range("E28").Value = -0.02493844809784762
$ws.Range("D29").Value = 0.0213160772790928
$ws.Range("E29").Value = -0.05073105892778018
$ws.Range("D30").Value = 0.02917765788465845
$ws.Range("E30").Value = -0.04124330117899244
$ws.Range("D31").Value = 0.01485367951534656
$ws.Range("E31").Value = -0.02921243281140451
$ws.Range("D32").Value = 0.01349396313038722
$ws.Range("E32").Value = -0.03311316600616454
$ws.Range("D33").Value = 0.01816632878302954
$ws.Range("E33").Value = -0.04809122459097692
$ws.Range("D34").Value = 0.04259487720749052
$ws.Range("E34").Value = -0.03075249112358269
$ws.Range("D35").Value = 0.01101057855094378
$ws.Range("E35").Value = -0.02522154055896397
$ws.Range("D36").Value = 0.01006025629178086
$ws.Range("E36").Value = 0.0006092784402471629
$ws.Range("D37").Value = 0.01044203640695333
$ws.Range("E37").Value = -0.03346810422282132
$ws.Range("D38").Value = 0.007452968303399573
$ws.Range("E38").Value = -0.003524672708962728
$ws.Range("D39").Value = 0.01209015189267614
$ws.Range("E39").Value = -0.0093817656964158
$ws.Range("D40").Value = 0.01751278677394744
$ws.Range("E40").Value = -0.02328222600794994
$ws.Range("D41").Value = 0.01732170907868992
$ws.Range("E41").Value = -0.01420142700329285
$ws.Range("D42").Value = 0.03211418744025269
$ws.Range("E42").Value = -0.02521180251241595
$ws.Range("D43").Value = 0.01146176396434747
$ws.Range("E43").Value = -0.02281198515696847
$ws.Range("D44").Value = 0.02185941342924077
$ws.Range("E44").Value = -0.02150966826326062
$ws.Range("D45").Value = 0.01238605650029097
$ws.Range("E45").Value = -0.04019815029427143
$ws.Range("D46").Value = 0.008602436677685217
$ws.Range("E46").Value = -0.03150765606595984
$ws.Range("D47").Value = 0.01347050842147672
$ws.Range("E47").Value = -0.008274132887588825
$ws.Range("D48").Value = 0.0108417046467882
$ws.Range("E48").Value = -0.02585669781931477
$ws.Range("D49").Value = 0.01583849583308051
$ws.Range("E49").Value = -0.02440469138727641
$ws.Range("D50").Value = 0.008557810184864843
$ws.Range("E50").Value = -0.02789339628502208
$ws.Range("D51").Value = 0.01183618430459328
$ws.Range("E51").Value = -0.07279116465863456
$ws.Range("D52").Value = 0.008297337824177469
$ws.Range("E52").Value = 0.001824212271973602
$ws.Range("D53").Value = 0.009977070257511633
$ws.Range("E53").Value = -0.004410215903106862
$ws.Range("D54").Value = 0.1358023893164359
$ws.Range("E54").Value = -0.0003942440370587885
$ws.Range("D55").Value = 0.0437965088543931
$ws.Range("E55").Value = -0.02330662782228698
$ws.Range("D56").Value = 1
$ws.Range("E56").Value = -0.02001555307391667

# Restore sheet protection (contents locked; matches original workbook)
$ws.Protect()
